# Update the statistical description table (rows 2-19, columns B-I)
# with the new values for the re-sampled dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 3038
$ws.Range("C2").Value2 = 476.0121790651745
$ws.Range("D2").Value2 = 50.37639214489838
$ws.Range("E2").Value2 = 407
$ws.Range("F2").Value2 = 441
$ws.Range("G2").Value2 = 453
$ws.Range("H2").Value2 = 512
$ws.Range("I2").Value2 = 672

$ws.Range("B3").Value2 = 8197
$ws.Range("C3").Value2 = 57.26076491399292
$ws.Range("D3").Value2 = 6.131069641327043
$ws.Range("E3").Value2 = 43.05
$ws.Range("F3").Value2 = 53.13
$ws.Range("G3").Value2 = 57.48
$ws.Range("H3").Value2 = 60.87
$ws.Range("I3").Value2 = 77.73

$ws.Range("B4").Value2 = 8197
$ws.Range("C4").Value2 = 11.10812736366964
$ws.Range("D4").Value2 = 14.77310032576502
$ws.Range("E4").Value2 = 0.85
$ws.Range("F4").Value2 = 7.26
$ws.Range("G4").Value2 = 11.82
$ws.Range("H4").Value2 = 13.44
$ws.Range("I4").Value2 = 402.27

$ws.Range("B5").Value2 = 8197
$ws.Range("C5").Value2 = 323.2362742466756
$ws.Range("D5").Value2 = 1.951333985013952
$ws.Range("E5").Value2 = 319.47
$ws.Range("F5").Value2 = 321.45
$ws.Range("G5").Value2 = 323.5
$ws.Range("H5").Value2 = 324.94
$ws.Range("I5").Value2 = 326.61

$ws.Range("B6").Value2 = 8197
$ws.Range("C6").Value2 = 26.0875613029157
$ws.Range("D6").Value2 = 1.650431260855432
$ws.Range("E6").Value2 = 0
$ws.Range("F6").Value2 = 25.48
$ws.Range("G6").Value2 = 26.35
$ws.Range("H6").Value2 = 27.19
$ws.Range("I6").Value2 = 29.99

$ws.Range("B7").Value2 = 8197
$ws.Range("C7").Value2 = -46.58728803220691
$ws.Range("D7").Value2 = 10.20621201600229
$ws.Range("E7").Value2 = -128
$ws.Range("F7").Value2 = -57
$ws.Range("G7").Value2 = -43
$ws.Range("H7").Value2 = -38
$ws.Range("I7").Value2 = -28

$ws.Range("B8").Value2 = 8196
$ws.Range("C8").Value2 = 10.26993655441679
$ws.Range("D8").Value2 = 1.771792829661422
$ws.Range("E8").Value2 = -23.5
$ws.Range("F8").Value2 = 9
$ws.Range("G8").Value2 = 10.2
$ws.Range("H8").Value2 = 11.8
$ws.Range("I8").Value2 = 15.5

$ws.Range("B9").Value2 = 8197
$ws.Range("C9").Value2 = 9.386726851287056
$ws.Range("D9").Value2 = 1.676426131130018
$ws.Range("E9").Value2 = 7
$ws.Range("F9").Value2 = 8
$ws.Range("G9").Value2 = 9
$ws.Range("H9").Value2 = 11
$ws.Range("I9").Value2 = 12

$ws.Range("B10").Value2 = 8197
$ws.Range("C10").Value2 = 867.8380261071123
$ws.Range("D10").Value2 = 0.46382526644809
$ws.Range("E10").Value2 = 867.1
$ws.Range("F10").Value2 = 867.5
$ws.Range("G10").Value2 = 867.9
$ws.Range("H10").Value2 = 868.3
$ws.Range("I10").Value2 = 868.5

$ws.Range("B11").Value2 = 8196
$ws.Range("C11").Value2 = 1632.619204489995
$ws.Range("D11").Value2 = 1050.750478570848
$ws.Range("E11").Value2 = 1
$ws.Range("F11").Value2 = 757
$ws.Range("G11").Value2 = 1474.5
$ws.Range("H11").Value2 = 2487.25
$ws.Range("I11").Value2 = 3898

$ws.Range("B12").Value2 = 8197
$ws.Range("C12").Value2 = 1743.127729657192
$ws.Range("D12").Value2 = 1129.67439147521
$ws.Range("E12").Value2 = 0
$ws.Range("F12").Value2 = 801
$ws.Range("G12").Value2 = 1576
$ws.Range("H12").Value2 = 2641
$ws.Range("I12").Value2 = 4224

$ws.Range("B13").Value2 = 8197
$ws.Range("C13").Value2 = 0.525505163352446
$ws.Range("D13").Value2 = 0.5409177378327277
$ws.Range("E13").Value2 = 0.061696
$ws.Range("F13").Value2 = 0.123392
$ws.Range("G13").Value2 = 0.246784
$ws.Range("H13").Value2 = 0.823296
$ws.Range("I13").Value2 = 1.974272

$ws.Range("B14").Value2 = 8197
$ws.Range("C14").Value2 = 0.04610345248261561
$ws.Range("D14").Value2 = 0.02173463530691788
$ws.Range("E14").Value2 = 0.02
$ws.Range("F14").Value2 = 0.02
$ws.Range("G14").Value2 = 0.05
$ws.Range("H14").Value2 = 0.07
$ws.Range("I14").Value2 = 0.07

$ws.Range("B15").Value2 = 8197
$ws.Range("C15").Value2 = 14
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 14
$ws.Range("F15").Value2 = 14
$ws.Range("G15").Value2 = 14
$ws.Range("H15").Value2 = 14
$ws.Range("I15").Value2 = 14

$ws.Range("B16").Value2 = 8197
$ws.Range("C16").Value2 = 1
$ws.Range("D16").Value2 = 0
$ws.Range("E16").Value2 = 1
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 1
$ws.Range("H16").Value2 = 1
$ws.Range("I16").Value2 = 1

$ws.Range("B17").Value2 = 8197
$ws.Range("C17").Value2 = 0.4000000000000001
$ws.Range("D17").Value2 = 0.00000000000005745754642241973
$ws.Range("E17").Value2 = 0.4
$ws.Range("F17").Value2 = 0.4
$ws.Range("G17").Value2 = 0.4
$ws.Range("H17").Value2 = 0.4
$ws.Range("I17").Value2 = 0.4

$ws.Range("B18").Value2 = 8197
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 0
$ws.Range("E18").Value2 = 3
$ws.Range("F18").Value2 = 3
$ws.Range("G18").Value2 = 3
$ws.Range("H18").Value2 = 3
$ws.Range("I18").Value2 = 3

$ws.Range("B19").Value2 = 8197
$ws.Range("C19").Value2 = 62.9872880322069
$ws.Range("D19").Value2 = 10.20621201600226
$ws.Range("E19").Value2 = 44.4
$ws.Range("F19").Value2 = 54.4
$ws.Range("G19").Value2 = 59.4
$ws.Range("H19").Value2 = 73.4
$ws.Range("I19").Value2 = 144.4
